# Re-colour the deck's theme (ppt/theme/theme1.xml) from the custom
# "Integral" palette back to the stock PowerPoint "Office Theme" palette.
#
# dk1 / lt1 (black / white) are unchanged between the two palettes, so
# only the remaining 10 theme colours (dk2, lt2, accent1-6, hlink,
# folHlink) need to be rewritten.

$p = $ppt.ActivePresentation
$tcs = $p.Designs.Item(1).SlideMaster.Theme.ThemeColorScheme

# Index -> (slot, new RGB as 0xBBGGRR OLE colour)
$tcs.Item(3).RGB  = 0x6A5444   # dk2      -> 44546A
$tcs.Item(4).RGB  = 0xE6E6E7   # lt2      -> E7E6E6
$tcs.Item(5).RGB  = 0xD59B5B   # accent1  -> 5B9BD5
$tcs.Item(6).RGB  = 0x317DED   # accent2  -> ED7D31
$tcs.Item(7).RGB  = 0xA5A5A5   # accent3  -> A5A5A5
$tcs.Item(8).RGB  = 0x00C0FF   # accent4  -> FFC000
$tcs.Item(9).RGB  = 0xC47244   # accent5  -> 4472C4
$tcs.Item(10).RGB = 0x47AD70   # accent6  -> 70AD47
$tcs.Item(11).RGB = 0xC16305   # hlink    -> 0563C1
$tcs.Item(12).RGB = 0x724F95   # folHlink -> 954F72
